# Apply crypto price/volume updates per commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.487.09"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "3.769.48"
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "615.73"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.01"
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").Value = "3.767.99"
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.527"
$ws.Range("E9").Value = "  -1.43%  "
$ws.Range("E10").Value = "  -2.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.69"
$ws.Range("E11").Value = "  +5.99%  "
$ws.Range("E12").Value = "  -1.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.02"
$ws.Range("E13").Value = "  -2.32%  "
$ws.Range("E14").Value = "  -3.43%  "
$ws.Range("D15").Value = "4.396.75"
$ws.Range("E15").Value = "  -0.31%  "
$ws.Range("D16").Value = "3.768.27"
$ws.Range("E16").Value = "  -0.42%  "
$ws.Range("D17").Value = "69.521.76"
$ws.Range("E17").Value = "  -0.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.52"
$ws.Range("E18").Value = "  -1.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.119"
$ws.Range("E19").Value = "  -3.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "508.00"
$ws.Range("E20").Value = "  -0.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.30"
$ws.Range("E21").Value = "  -3.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.37"
$ws.Range("E22").Value = "  -1.23%  "
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.51"
$ws.Range("E24").Value = "  +0.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.35"
$ws.Range("E25").Value = "  -1.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.80"
$ws.Range("E26").Value = "  -2.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000136"
$ws.Range("E27").Value = "  -2.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.54"
$ws.Range("E28").Value = "  -4.03%  "
$ws.Range("E29").Value = "  +0.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.50"
$ws.Range("E30").Value = "  +0.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.95"
$ws.Range("E31").Value = "  +3.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.98"
$ws.Range("E32").Value = "  +2.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.67"
$ws.Range("E33").Value = "  -1.97%  "
$ws.Range("E34").Value = "  -1.32%  "
$ws.Range("E35").Value = "  -0.14%  "
$ws.Range("E36").Value = "  -1.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.13"
$ws.Range("E37").Value = "  -1.17%  "
$ws.Range("E38").Value = "  +3.53%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.339"
$ws.Range("E39").Value = "  +2.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "451.50"
$ws.Range("E40").Value = "  +8.02%  "
$ws.Range("E41").Value = "  -2.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "49.88"
$ws.Range("E42").Value = "  -2.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.98"
$ws.Range("E43").Value = "  +5.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "44.50"
$ws.Range("E44").Value = "  -0.99%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.55"
$ws.Range("E45").Value = "  -2.16%  "
$ws.Range("D46").Value = "2.955.64"
$ws.Range("E46").Value = "  -2.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0358"
$ws.Range("E47").Value = "  -1.40%  "
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "139.14"
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "27.16"
$ws.Range("E50").Value = "  -0.92%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.45"
$ws.Range("E51").Value = "  -1.04%  "
